$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "JSU(-0.9948501779716727, 1.0695082297450007, 0.34650546570203294, 2.4871216326779506)"
$ws.Range("C2").Value = "NIG(1.1832031182873146, 0.8826769137504804, 4.061243333425183, 6.0466985364330395)"
$ws.Range("D2").Value = "NIG(1.6160096479263695, 1.2161136324198742, -0.01652738226860799, 4.119951486067391)"
$ws.Range("E2").Value = "NIG(1.1867929129822894, 0.9473461574450233, 4.201551075635541, 5.001887420949407)"
